$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column numeric-looking values are prefixed with a leading apostrophe so
# Excel stores them as text (matching the source inlineStr cells) instead of
# silently converting them to numbers (e.g. "1.000" -> 1).

$ws.Range('D2').Value = '30.516.96'
$ws.Range('E2').Value = '  +0.34%  '

$ws.Range('D3').Value = '1.915.37'
$ws.Range('E3').Value = '  -0.04%  '

$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').Value = '''245.54'
$ws.Range('E5').Value = '  +1.28%  '

$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.16%  '

$ws.Range('D7').Value = '''0.4822'
$ws.Range('E7').Value = '  +2.85%  '

$ws.Range('D8').Value = '''0.2890'
$ws.Range('E8').Value = '  +1.17%  '

$ws.Range('D9').Value = '''0.06722'
$ws.Range('E9').Value = '  -1.40%  '

$ws.Range('D10').Value = '''110.40'
$ws.Range('E10').Value = '  -0.71%  '

$ws.Range('D11').Value = '''19.26'
$ws.Range('E11').Value = '  +4.71%  '

$ws.Range('D12').Value = '1.917.69'
$ws.Range('E12').Value = '  +0.22%  '

$ws.Range('D13').Value = '''0.07564'
$ws.Range('E13').Value = '  -2.07%  '

$ws.Range('D14').Value = '''5.253'
$ws.Range('E14').Value = '  -0.12%  '

$ws.Range('D15').Value = '''0.6699'
$ws.Range('E15').Value = '  +1.72%  '

$ws.Range('D16').Value = '''288.42'
$ws.Range('E16').Value = '  -3.25%  '

$ws.Range('D17').Value = '30.516.70'
$ws.Range('E17').Value = '  +0.31%  '

$ws.Range('D18').Value = '''0.000007603'
$ws.Range('E18').Value = '  -0.23%  '

$ws.Range('D19').Value = '''1.000'
$ws.Range('E19').Value = '  +0.05%  '

$ws.Range('D20').Value = '''12.90'
$ws.Range('E20').Value = '  -0.13%  '

$ws.Range('D21').Value = '2.165.89'
$ws.Range('E21').Value = '  +1.14%  '

$ws.Range('D22').Value = '''5.485'
$ws.Range('E22').Value = '  +4.69%  '

$ws.Range('E23').Value = '  +0.32%  '

$ws.Range('D24').Value = '''6.394'
$ws.Range('E24').Value = '  +2.83%  '

$ws.Range('D25').Value = '''9.460'
$ws.Range('E25').Value = '  +1.65%  '

$ws.Range('D26').Value = '''164.56'
$ws.Range('E26').Value = '  -2.11%  '

$ws.Range('D27').Value = '''20.34'
$ws.Range('E27').Value = '  -6.00%  '

$ws.Range('D28').Value = '''2.138'
$ws.Range('E28').Value = '  +2.56%  '

$ws.Range('D29').Value = '''0.1062'
$ws.Range('E29').Value = '  -0.80%  '

$ws.Range('D30').Value = '''1.408'

$ws.Range('D31').Value = '''4.158'
$ws.Range('E31').Value = '  -0.24%  '

$ws.Range('D32').Value = '''4.028'
$ws.Range('E32').Value = '  +1.09%  '

$ws.Range('E33').Value = '  -1.43%  '

$ws.Range('D34').Value = '''0.7287'
$ws.Range('E34').Value = '  -1.38%  '

$ws.Range('D35').Value = '''1.134'
$ws.Range('E35').Value = '  -1.65%  '

$ws.Range('E36').Value = '  -1.33%  '

$ws.Range('D37').Value = '''0.9993'
$ws.Range('E37').Value = '  +0.15%  '

$ws.Range('D38').Value = '''2.737'
$ws.Range('E38').Value = '  -0.20%  '

$ws.Range('D39').Value = '''2.672'
$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('D40').Value = '''110.68'
$ws.Range('E40').Value = '  +0.54%  '

$ws.Range('D41').Value = '''2.012'
$ws.Range('E41').Value = '  -2.22%  '

$ws.Range('D42').Value = '''0.4428'
$ws.Range('E42').Value = '  +3.85%  '

$ws.Range('D43').Value = '''0.8638'
$ws.Range('E43').Value = '  -1.04%  '

$ws.Range('E44').Value = '  +1.35%  '

$ws.Range('D45').Value = '''1.000'
$ws.Range('E45').Value = '  +0.13%  '

$ws.Range('D46').Value = '''68.14'
$ws.Range('E46').Value = '  +0.94%  '

$ws.Range('D47').Value = '''7.344'
$ws.Range('E47').Value = '  +2.30%  '

$ws.Range('D48').Value = '''48.89'
$ws.Range('E48').Value = '  -6.03%  '

$ws.Range('D49').Value = '''9.264'
$ws.Range('E49').Value = '  +0.22%  '

$ws.Range('D50').Value = '''0.1242'
$ws.Range('E50').Value = '  +2.62%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '''0.2512'
$ws.Range('E51').Value = '  +3.85%  '
